# Add new "WhatsApp B" submission-time rows to the bottom of each of the
# four tracking sheets, matching the target diff.

$wb = $excel.ActiveWorkbook

function Add-Row {
    param($ws, $row, $aVal, $bVal, $cVal, $dVal, $eVal)

    $rng = $ws.Range("A" + $row + ":E" + $row)
    $rng.Style = "Normal"

    $ws.Cells.Item($row, 1).Value = $aVal
    $ws.Cells.Item($row, 2).Value = $bVal
    $ws.Cells.Item($row, 3).Value = $cVal
    $ws.Cells.Item($row, 4).Value = $dVal
    $ws.Cells.Item($row, 5).Value = $eVal
}

# Sheet 1: "Submit orders" -> add row 82
$ws1 = $wb.Worksheets.Item("Submit orders")
Add-Row $ws1 82 "10.13.2022 12:51 (Kyiv+Israel) 09:51 (UTC) 18:51 (Japan) 15:21 (India)" "***" "***" 1.625 -0.4870000000000001

# Sheet 2: "Submit internet survey" -> add row 77
$ws2 = $wb.Worksheets.Item("Submit internet survey")
Add-Row $ws2 77 "10.13.2022 12:40 (Kyiv+Israel) 09:40 (UTC) 18:40 (Japan) 15:10 (India)" "***" "***" 0.967 -0.202

# Sheet 3: "Submit a phone survey" -> add rows 71 and 72
$ws3 = $wb.Worksheets.Item("Submit a phone survey")
Add-Row $ws3 71 "10.13.2022 12:53 (Kyiv+Israel) 09:53 (UTC) 18:53 (Japan) 15:23 (India)" "***" "***" 2.5 -0.9159999999999999
Add-Row $ws3 72 "10.13.2022 14:43 (Kyiv+Israel) 11:43 (UTC) 20:43 (Japan) 17:13 (India)" "***" "***" 2.041 -0.4569999999999999

# Sheet 4: "Checkertificate" -> add rows 79, 80, 81, 82
$ws4 = $wb.Worksheets.Item("Checkertificate")
Add-Row $ws4 79 "10.12.2022 10:48 (Kyiv+Israel) 07:48 (UTC) 16:48 (Japan) 13:18 (India)" 0.973 -0.3079999999999999 "***" "***"
Add-Row $ws4 80 "10.12.2022 10:50 (Kyiv+Israel) 07:50 (UTC) 16:50 (Japan) 13:20 (India)" "***" "***" 0.979 -0.05699999999999994
Add-Row $ws4 81 "10.13.2022 13:14 (Kyiv+Israel) 10:14 (UTC) 19:14 (Japan) 15:44 (India)" "***" "***" 1.081 -0.1589999999999999
Add-Row $ws4 82 "10.13.2022 14:40 (Kyiv+Israel) 11:40 (UTC) 20:40 (Japan) 17:10 (India)" "***" "***" 1.222 -0.2999999999999999
